# Updated symbol list on Fri Dec 23 08:30:30 UTC 2022 with GitHub Actions
#
# This script applies the cell-level price/label updates captured in the
# commit's OOXML diff. All affected cells are plain text cells
# (t="inlineStr" in the original workbook), including the "Price" column
# (D) whose contents look numeric (e.g. "245.99", "0.00002101"). To keep
# those as text (not auto-converted to a float, which would silently
# drop significant trailing/leading zeros such as "5.420" or
# "0.00002100"), we prefix them with a leading single-quote, exactly as
# a user typing into Excel would to force text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-9: price-only refresh
$ws.Range("D2").Value = "'245.99"
$ws.Range("D3").Value = "'21.98"
$ws.Range("D4").Value = "'5.373"
$ws.Range("D5").Value = "'0.05796"
$ws.Range("D6").Value = "'3.375"
$ws.Range("D7").Value = "'6.326"
$ws.Range("D8").Value = "'0.8089"
$ws.Range("D9").Value = "'0.9787"

# Rows 10-18: coin list reshuffled (ONE dropped to the bottom, others
# shifted up) with refreshed price/volume-label values
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1427"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07515"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03186"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03028"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'4.155"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09403"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001609"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04808"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005902"
$ws.Range("E18").Value = "17OneONE"

# Remaining scattered price / label refreshes
$ws.Range("D19").Value = "'0.005399"
$ws.Range("D20").Value = "'0.004091"
$ws.Range("D21").Value = "'0.0009963"
$ws.Range("D24").Value = "'2.244"
$ws.Range("D25").Value = "'0.3231"
$ws.Range("D26").Value = "'0.1297"
$ws.Range("D40").Value = "'0.03888"
$ws.Range("D41").Value = "'0.006342"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("D44").Value = "'0.006684"
$ws.Range("D48").Value = "'0.1469"
$ws.Range("D49").Value = "'0.00002101"

Write-Host "Applied 57 cell updates"
